# "Add analyzer to assets" -- tag a batch of Accounting-Title rows in
# Sheet1 with an Analyzer/ID value in column A (T004 / T005 / placeholder
# strings), matching the shared-string additions introduced by the
# target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Introduce the brand-new shared strings first, in the same order
# they first appear in the target workbook, so the rebuilt shared
# string table grows in a matching order as much as possible.
$ws.Range("A39").Value = "T004"
$ws.Range("A36").Value = "ts无此科目"
$ws.Range("A50").Value = "难以自动判定"
$ws.Range("A34").Value = "T004 - 仅固定资产"
$ws.Range("A60").Value = "T005"

# --- Remaining rows in the "资产负债表" block (rows 33-44) that reuse the
# already-existing "ts无此数据" placeholder or the strings just added above.
$ws.Range("A33").Value = "ts无此数据"
$ws.Range("A37").Value = "ts无此数据"
$ws.Range("A41").Value = "T004"
$ws.Range("A43").Value = "ts无此数据"
$ws.Range("A44").Value = "ts无此数据"

# --- "其它" non-core-asset block (rows 50-52)
$ws.Range("A51").Value = "T004 - 仅固定资产"
$ws.Range("A52").Value = "难以自动判定"

# --- 利润表 block (rows 60-74) all tagged T005
$ws.Range("A61").Value = "T005"
$ws.Range("A65").Value = "T005"
$ws.Range("A67").Value = "T005"
$ws.Range("A69").Value = "T005"
$ws.Range("A71").Value = "T005"
$ws.Range("A72").Value = "T005"
$ws.Range("A73").Value = "T005"
$ws.Range("A74").Value = "T005"

# Leave the selection near where the author ended up editing.
$ws.Range("A74").Select()
